$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "2025-03-05 17:28:24.329449"
$ws.Range("B2").Value = "fd31:1623:3a00:148:3554:8d02:4013:d1aa"
$ws.Range("D2").Value = "github.com."
$ws.Range("E2").Value = 1
